$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.036036036036036
$ws.Range("C2").Value = 0.000750750750750751
$ws.Range("D2").Value = 0.0015015015015015
$ws.Range("E2").Value = 0.003003003003003
$ws.Range("F2").Value = 0.00225225225225225
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0.990990990990991
$ws.Range("I2").Value = 0.021021021021021
$ws.Range("J2").Value = 0.990990990990991
$ws.Range("K2").Value = 0.0195195195195195
$ws.Range("L2").Value = 0.00600600600600601
$ws.Range("M2").Value = 0.021021021021021
$ws.Range("N2").Value = 0.000750750750750751
$ws.Range("O2").Value = 0.0315315315315315
$ws.Range("P2").Value = 0.0015015015015015
$ws.Range("Q2").Value = 0.99024024024024
$ws.Range("R2").Value = 0.00525525525525526
$ws.Range("S2").Value = 0.989489489489489
$ws.Range("T2").Value = 0.048048048048048
$ws.Range("U2").Value = 0.990990990990991
$ws.Range("V2").Value = 0.989489489489489
$ws.Range("W2").Value = 0.00225225225225225
$ws.Range("X2").Value = 0.00225225225225225

$ws.Range("B3").Value = 0.003003003003003
$ws.Range("C3").Value = 0.99024024024024
$ws.Range("D3").Value = 0.992492492492492
$ws.Range("E3").Value = 0.994744744744745
$ws.Range("F3").Value = 0.00225225225225225
$ws.Range("G3").Value = 0.996996996996997
$ws.Range("H3").Value = 0.000750750750750751
$ws.Range("I3").Value = 0.00225225225225225
$ws.Range("J3").Value = 0.00375375375375375
$ws.Range("K3").Value = 0.00225225225225225
$ws.Range("L3").Value = 0.018018018018018
$ws.Range("M3").Value = 0.972972972972973
$ws.Range("N3").Value = 0.00225225225225225
$ws.Range("O3").Value = 0.00375375375375375
$ws.Range("P3").Value = 0.996996996996997
$ws.Range("Q3").Value = 0.00600600600600601
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0.003003003003003
$ws.Range("T3").Value = 0.0015015015015015
$ws.Range("U3").Value = 0.00225225225225225
$ws.Range("V3").Value = 0.00600600600600601
$ws.Range("W3").Value = 0.00525525525525526
$ws.Range("X3").Value = 0.984234234234234

$ws.Range("B4").Value = 0.952702702702703
$ws.Range("C4").Value = 0.00375375375375375
$ws.Range("D4").Value = 0.000750750750750751
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0.000750750750750751
$ws.Range("G4").Value = 0.000750750750750751
$ws.Range("H4").Value = 0.0045045045045045
$ws.Range("I4").Value = 0.972222222222222
$ws.Range("J4").Value = 0.00375375375375375
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0.972222222222222
$ws.Range("M4").Value = 0.003003003003003
$ws.Range("N4").Value = 0.0045045045045045
$ws.Range("O4").Value = 0.963963963963964
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0.003003003003003
$ws.Range("S4").Value = 0.00525525525525526
$ws.Range("T4").Value = 0.00225225225225225
$ws.Range("U4").Value = 0.0015015015015015
$ws.Range("V4").Value = 0.003003003003003
$ws.Range("W4").Value = 0.987237237237237
$ws.Range("X4").Value = 0.00825825825825826

$ws.Range("B5").Value = 0.00600600600600601
$ws.Range("C5").Value = 0.00525525525525526
$ws.Range("D5").Value = 0.00525525525525526
$ws.Range("E5").Value = 0.00225225225225225
$ws.Range("F5").Value = 0.994744744744745
$ws.Range("G5").Value = 0.00225225225225225
$ws.Range("H5").Value = 0.00375375375375375
$ws.Range("I5").Value = 0.0045045045045045
$ws.Range("J5").Value = 0.0015015015015015
$ws.Range("K5").Value = 0.978228228228228
$ws.Range("L5").Value = 0.00375375375375375
$ws.Range("M5").Value = 0.003003003003003
$ws.Range("N5").Value = 0.992492492492492
$ws.Range("O5").Value = 0.000750750750750751
$ws.Range("P5").Value = 0.0015015015015015
$ws.Range("Q5").Value = 0.00375375375375375
$ws.Range("R5").Value = 0.991741741741742
$ws.Range("S5").Value = 0.00225225225225225
$ws.Range("T5").Value = 0.948198198198198
$ws.Range("U5").Value = 0.00525525525525526
$ws.Range("V5").Value = 0.0015015015015015
$ws.Range("W5").Value = 0.00525525525525526
$ws.Range("X5").Value = 0.00525525525525526

